# Split " años" (the measurement-unit suffix after the age number) into its
# own run and change it to "meses", while leaving the preceding ": " and "63"
# runs untouched.
#
# Strategy: a direct Range.Text assignment that spans a run boundary causes
# this engine to coalesce the edited run with its identically-formatted
# neighbour, which would merge "63" and " años" together. Toggling a
# character-formatting property (Bold) on a sub-range first forces the run
# to split without merging anything, isolating the "años" text in its own
# run. We can then safely replace its text and undo the formatting toggle.

$d = $word.ActiveDocument

# Locate the table cell that holds "Edad del paciente: 63 años".
$ageCell = $null
$t = $d.Tables.Item(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $labelCell = $t.Cell($r, 1)
    if ($labelCell.Range.Text.StartsWith("Edad del paciente")) {
        $ageCell = $t.Cell($r, 2)
        break
    }
}

$cellRange = $ageCell.Range

# Find the exact extent of the " años" text inside that cell.
$unitRange = $cellRange.Duplicate
$unitRange.Find.Execute(" años", $false, $false, $false, $false, $false, $true, 0, $false, "", 0)

# Isolate just the "años" word (without the leading space) into its own run
# by toggling Bold on/off -- a pure formatting op splits runs but never
# merges them, unlike a Text assignment.
$wordRange = $d.Range($unitRange.Start + 1, $unitRange.End)
$wordRange.Bold = 1
$wordRange.Text = "meses"
$wordRange.Bold = 0
